$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-14 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-15 Thursday", 2) | Out-Null
$d.Content.Find.Execute("970÷2=485, 0", $true, $false, $false, $false, $false, $true, 1, $false, "555÷9=61, 6", 2) | Out-Null
$d.Content.Find.Execute("627÷7=89, 4", $true, $false, $false, $false, $false, $true, 1, $false, "666÷3=222, 0", 2) | Out-Null
$d.Content.Find.Execute("106÷9=11, 7", $true, $false, $false, $false, $false, $true, 1, $false, "324÷7=46, 2", 2) | Out-Null
$d.Content.Find.Execute("939÷4=234, 3", $true, $false, $false, $false, $false, $true, 1, $false, "474÷3=158, 0", 2) | Out-Null
$d.Content.Find.Execute("580÷3=193, 1", $true, $false, $false, $false, $false, $true, 1, $false, "461÷5=92, 1", 2) | Out-Null
$d.Content.Find.Execute("687÷7=98, 1", $true, $false, $false, $false, $false, $true, 1, $false, "196÷9=21, 7", 2) | Out-Null
$d.Content.Find.Execute("244÷3=81, 1", $true, $false, $false, $false, $false, $true, 1, $false, "259÷5=51, 4", 2) | Out-Null
$d.Content.Find.Execute("895÷9=99, 4", $true, $false, $false, $false, $false, $true, 1, $false, "145÷4=36, 1", 2) | Out-Null
$d.Content.Find.Execute("341÷6=56, 5", $true, $false, $false, $false, $false, $true, 1, $false, "680÷3=226, 2", 2) | Out-Null
$d.Content.Find.Execute("978÷6=163, 0", $true, $false, $false, $false, $false, $true, 1, $false, "976÷3=325, 1", 2) | Out-Null
$d.Content.Find.Execute("611÷4=152, 3", $true, $false, $false, $false, $false, $true, 1, $false, "889÷5=177, 4", 2) | Out-Null
$d.Content.Find.Execute("771÷3=257, 0", $true, $false, $false, $false, $false, $true, 1, $false, "870÷9=96, 6", 2) | Out-Null
$d.Content.Find.Execute("815÷3=271, 2", $true, $false, $false, $false, $false, $true, 1, $false, "403÷4=100, 3", 2) | Out-Null
$d.Content.Find.Execute("281÷9=31, 2", $true, $false, $false, $false, $false, $true, 1, $false, "158÷6=26, 2", 2) | Out-Null
$d.Content.Find.Execute("568÷9=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "155÷8=19, 3", 2) | Out-Null
$d.Content.Find.Execute("557÷9=61, 8", $true, $false, $false, $false, $false, $true, 1, $false, "198÷9=22, 0", 2) | Out-Null
$d.Content.Find.Execute("158÷7=22, 4", $true, $false, $false, $false, $false, $true, 1, $false, "877÷3=292, 1", 2) | Out-Null
$d.Content.Find.Execute("163÷3=54, 1", $true, $false, $false, $false, $false, $true, 1, $false, "309÷8=38, 5", 2) | Out-Null
$d.Content.Find.Execute("251÷7=35, 6", $true, $false, $false, $false, $false, $true, 1, $false, "166÷3=55, 1", 2) | Out-Null
$d.Content.Find.Execute("555÷8=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "964÷7=137, 5", 2) | Out-Null
$d.Content.Find.Execute("527÷3=175, 2", $true, $false, $false, $false, $false, $true, 1, $false, "133÷3=44, 1", 2) | Out-Null
$d.Content.Find.Execute("988÷6=164, 4", $true, $false, $false, $false, $false, $true, 1, $false, "540÷7=77, 1", 2) | Out-Null
$d.Content.Find.Execute("250÷4=62, 2", $true, $false, $false, $false, $false, $true, 1, $false, "945÷8=118, 1", 2) | Out-Null
$d.Content.Find.Execute("142÷4=35, 2", $true, $false, $false, $false, $false, $true, 1, $false, "497÷9=55, 2", 2) | Out-Null
$d.Content.Find.Execute("879÷6=146, 3", $true, $false, $false, $false, $false, $true, 1, $false, "124÷8=15, 4", 2) | Out-Null
